$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Formula = "=33/50"
$ws.Range("H8").Value = "need a 67 on the final"
$ws.Range("H9").Value = "need a 73 on programming assignment #5"

$ws.Range("B9").Select()
